$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-28 Monday" "2025-07-29 Tuesday"

Replace-Text "18×55=" "18×38="
Replace-Text "58×63=" "53×35="
Replace-Text "35×77=" "21×68="
Replace-Text "39×83=" "90×38="
Replace-Text "28×47=" "76×18="
Replace-Text "33×77=" "27×41="
Replace-Text "42×64=" "31×20="
Replace-Text "13×87=" "92×98="
Replace-Text "29×79=" "69×54="
Replace-Text "31×79=" "27×96="
Replace-Text "14×49=" "82×70="
Replace-Text "20×87=" "73×13="
Replace-Text "24×96=" "60×68="
Replace-Text "53×26=" "86×11="
Replace-Text "33×56=" "58×26="
Replace-Text "49×95=" "59×55="
Replace-Text "56×96=" "80×75="
Replace-Text "15×80=" "80×62="
Replace-Text "76×19=" "77×92="
Replace-Text "40×56=" "50×48="
Replace-Text "28×76=" "72×88="
Replace-Text "34×84=" "53×14="
Replace-Text "60×78=" "24×59="
Replace-Text "81×89=" "94×29="
Replace-Text "48×55=" "95×62="
